# Auto-generated edit script: update TPM-derived metrics per commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5053303333333333
$ws.Range("H2").Value = 1.515991
$ws.Range("I2").Value = 0.006927186824079787
$ws.Range("J2").Value = 0.006927186824079787
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 0.9696815770587777
$ws.Range("R2").Value = 8.727134193529
$ws.Range("S2").Value = 0.00004519309958257517
$ws.Range("T2").Value = 0.00004519309958257517
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5053303333333333
$ws.Range("H3").Value = 1.515991
$ws.Range("I3").Value = 0.006927186824079787
$ws.Range("J3").Value = 0.006927186824079787
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 91.61052322193711
$ws.Range("R3").Value = 824.494708997434
$ws.Range("S3").Value = 0.004269611382469177
$ws.Range("T3").Value = 0.004269611382469177
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5053303333333333
$ws.Range("H4").Value = 1.515991
$ws.Range("I4").Value = 0.006927186824079787
$ws.Range("J4").Value = 0.006927186824079787
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 14.67866994276078
$ws.Range("R4").Value = 132.108029484847
$ws.Range("S4").Value = 0.0006841159078994547
$ws.Range("T4").Value = 0.0006841159078994548
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5053303333333333
$ws.Range("H5").Value = 1.515991
$ws.Range("I5").Value = 0.006927186824079787
$ws.Range("J5").Value = 0.006927186824079787
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 41.37367107159511
$ws.Range("R5").Value = 372.363039644356
$ws.Range("S5").Value = 0.00192826643412858
$ws.Range("T5").Value = 0.00192826643412858
# Row 6
$ws.Range("G6").Value = 62.87584200000001
$ws.Range("I6").Value = 0.8619168008028857
$ws.Range("J6").Value = 0.8619168008028857
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 120.652851427466
$ws.Range("R6").Value = 1085.875662847194
$ws.Range("S6").Value = 0.005623161724926327
$ws.Range("T6").Value = 0.005623161724926327
# Row 7
$ws.Range("G7").Value = 62.87584200000001
$ws.Range("I7").Value = 0.8619168008028857
$ws.Range("J7").Value = 0.8619168008028857
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.5312473702394016
$ws.Range("T7").Value = 0.5312473702394016
# Row 8
$ws.Range("G8").Value = 62.87584200000001
$ws.Range("I8").Value = 0.8619168008028857
$ws.Range("J8").Value = 0.8619168008028857
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 1826.396856098438
$ws.Range("R8").Value = 16437.57170488594
$ws.Range("S8").Value = 0.08512127789961682
$ws.Range("T8").Value = 0.08512127789961683
# Row 9
$ws.Range("G9").Value = 62.87584200000001
$ws.Range("I9").Value = 0.8619168008028857
$ws.Range("J9").Value = 0.8619168008028857
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 5147.928461166824
$ws.Range("R9").Value = 46331.35615050142
$ws.Range("S9").Value = 0.239924990938941
$ws.Range("T9").Value = 0.239924990938941
# Row 10
$ws.Range("G10").Value = 4.893887666666667
$ws.Range("H10").Value = 14.681663
$ws.Range("I10").Value = 0.0670865608629469
$ws.Range("J10").Value = 0.06708656086294688
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 9.390912038188555
$ws.Range("R10").Value = 84.518208343697
$ws.Range("S10").Value = 0.0004376740086166801
$ws.Range("T10").Value = 0.00043767400861668
# Row 11
$ws.Range("G11").Value = 4.893887666666667
$ws.Range("H11").Value = 14.681663
$ws.Range("I11").Value = 0.0670865608629469
$ws.Range("J11").Value = 0.06708656086294688
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 887.2050224560404
$ws.Range("R11").Value = 7984.845202104362
$ws.Range("S11").Value = 0.04134918707193946
$ws.Range("T11").Value = 0.04134918707193945
# Row 12
$ws.Range("G12").Value = 4.893887666666667
$ws.Range("H12").Value = 14.681663
$ws.Range("I12").Value = 0.0670865608629469
$ws.Range("J12").Value = 0.06708656086294688
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 142.1560453774746
$ws.Range("R12").Value = 1279.404408397271
$ws.Range("S12").Value = 0.006625342243271124
$ws.Range("T12").Value = 0.006625342243271123
# Row 13
$ws.Range("G13").Value = 4.893887666666667
$ws.Range("H13").Value = 14.681663
$ws.Range("I13").Value = 0.0670865608629469
$ws.Range("J13").Value = 0.06708656086294688
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 400.6846318652342
$ws.Range("R13").Value = 3606.161686787108
$ws.Range("S13").Value = 0.01867435753911964
$ws.Range("T13").Value = 0.01867435753911963
# Row 14
$ws.Range("G14").Value = 4.673793
$ws.Range("H14").Value = 14.021379
$ws.Range("I14").Value = 0.06406945151008747
$ws.Range("J14").Value = 0.06406945151008747
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 8.968570988388999
$ws.Range("R14").Value = 80.71713889550099
$ws.Range("S14").Value = 0.0004179903293832405
$ws.Range("T14").Value = 0.0004179903293832405
# Row 15
$ws.Range("G15").Value = 4.673793
$ws.Range("H15").Value = 14.021379
$ws.Range("I15").Value = 0.06406945151008747
$ws.Range("J15").Value = 0.06406945151008747
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 847.3044143949941
$ws.Range("R15").Value = 7625.739729554946
$ws.Range("S15").Value = 0.03948957439477825
$ws.Range("T15").Value = 0.03948957439477825
# Row 16
$ws.Range("G16").Value = 4.673793
$ws.Range("H16").Value = 14.021379
$ws.Range("I16").Value = 0.06406945151008747
$ws.Range("J16").Value = 0.06406945151008747
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 135.762807617827
$ws.Range("R16").Value = 1221.865268560443
$ws.Range("S16").Value = 0.006327378213054925
$ws.Range("T16").Value = 0.006327378213054926
# Row 17
$ws.Range("G17").Value = 4.673793
$ws.Range("H17").Value = 14.021379
$ws.Range("I17").Value = 0.06406945151008747
$ws.Range("J17").Value = 0.06406945151008747
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 382.664489905396
$ws.Range("R17").Value = 3443.980409148564
$ws.Range("S17").Value = 0.01783450857287105
$ws.Range("T17").Value = 0.01783450857287105

